{"js": "// Replace the quoted phrase \"5 Bars\" with \"Pack of 4\" in the\n// \"Use \"5 Bars\" while ordering Soaps\" paragraph.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"5 Bars\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the text \"5 Bars\" to replace.');\n}\n\n// There should be exactly one occurrence; replace its text in place so the\n// surrounding quotes / sentence (\"Use \"\u2026\" while ordering Soaps\") are kept.\nsearchResults.items[0].insertText(\"Pack of 4\", \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the quoted phrase \"5 Bars\" with \"Pack of 4\" in the\n# \"Use \"5 Bars\" while ordering Soaps\" paragraph.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"5 Bars\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Pack of 4\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
